$wb = $excel.ActiveWorkbook

# --- Sheet "ger-pol-fin-research" ---
$ws1 = $wb.Worksheets.Item("ger-pol-fin-research")

# Insert 5 new rows (source data for Vestas' total capacity, by year) right after row 28
# Do this first so the formulas below can reference the final row numbers directly.
$ws1.Rows("29:33").Insert()

# Add new data point for Vestas/Germany 2021 (row 22): total capacity + its % of market
$ws1.Range("G22").Value = 16438
$ws1.Range("I22").Formula = "=G22/G39"

# Add new data point for Vestas/Poland 2021 (row 28): total capacity + its % of market
$ws1.Range("G28").Value = 3559
$ws1.Range("I28").Formula = "=G28/G45"

$ws1.Range("A29").Value = 2017
$ws1.Range("B29").Value = "Vestas"
$ws1.Range("C29").Value = "Total"
$ws1.Range("E29").Value = 8779

$ws1.Range("A30").Value = 2018
$ws1.Range("B30").Value = "Vestas"
$ws1.Range("C30").Value = "Total"
$ws1.Range("E30").Value = 10847

$ws1.Range("A31").Value = 2019
$ws1.Range("B31").Value = "Vestas"
$ws1.Range("C31").Value = "Total"
$ws1.Range("E31").Value = 12884

$ws1.Range("A32").Value = 2020
$ws1.Range("B32").Value = "Vestas"
$ws1.Range("C32").Value = "Total"
$ws1.Range("E32").Value = 17212

$ws1.Range("A33").Value = 2021
$ws1.Range("B33").Value = "Vestas"
$ws1.Range("C33").Value = "Total"
$ws1.Range("E33").Value = 16594

# --- Sheet "Sheet1" ---
$ws2 = $wb.Worksheets.Item("Sheet1")

# Insert 10 new rows mirroring the Vestas % of market figures
$ws2.Rows("12:21").Insert()

$ws2.Range("A12").Value = 2017
$ws2.Range("B12").Value = "Vestas"
$ws2.Range("C12").Value = "Germany"
$ws2.Range("D12").Value = 0.152

$ws2.Range("A13").Value = 2018
$ws2.Range("B13").Value = "Vestas"
$ws2.Range("C13").Value = "Germany"
$ws2.Range("D13").Value = 0.123

$ws2.Range("A14").Value = 2019
$ws2.Range("B14").Value = "Vestas"
$ws2.Range("C14").Value = "Germany"
$ws2.Range("D14").Value = 0.035

$ws2.Range("A15").Value = 2020
$ws2.Range("B15").Value = "Vestas"
$ws2.Range("C15").Value = "Germany"
$ws2.Range("D15").Value = 0.029

$ws2.Range("A16").Value = 2021
$ws2.Range("B16").Value = "Vestas"
$ws2.Range("C16").Value = "Germany"
$ws2.Range("D16").Value = 0.036
$ws2.Range("E16").Value = 0.25668332292317303

$ws2.Range("A17").Value = 2017
$ws2.Range("B17").Value = "Vestas"
$ws2.Range("C17").Value = "Poland"
$ws2.Range("D17").Value = 0

$ws2.Range("A18").Value = 2018
$ws2.Range("B18").Value = "Vestas"
$ws2.Range("C18").Value = "Poland"
$ws2.Range("D18").Value = 0

$ws2.Range("A19").Value = 2019
$ws2.Range("B19").Value = "Vestas"
$ws2.Range("C19").Value = "Poland"
$ws2.Range("D19").Value = 0.006

$ws2.Range("A20").Value = 2020
$ws2.Range("B20").Value = "Vestas"
$ws2.Range("C20").Value = "Poland"
$ws2.Range("D20").Value = 0.024

$ws2.Range("A21").Value = 2021
$ws2.Range("B21").Value = "Vestas"
$ws2.Range("C21").Value = "Poland"
$ws2.Range("D21").Value = 0.045
$ws2.Range("E21").Value = 0.50014052838673417
